# Update crypto price/volume snapshot (GitHub Actions scheduled refresh).
# Note: several "Price" values are plain decimal numbers (e.g. 18.48); a
# leading apostrophe is used so Excel stores them as text (quote-prefix),
# matching the source data which keeps every Price/Volume cell as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.013.60'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.638.28'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').Value = '''18.48'
$ws.Range('E10').Value = '  -6.10%  '
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.21'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.636.11'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').Value = '26.024.23'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''61.94'
$ws.Range('E17').Value = '  -2.37%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0744'
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '''192.51'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Value = '''4.25'
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').Value = '''9.76'
$ws.Range('D23').Value = '''6.11'
$ws.Range('E23').Value = '  -2.23%  '
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''143.93'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''1.78'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '''6.84'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = '''15.25'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  -2.59%  '
$ws.Range('D32').Value = '''3.15'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('E34').Value = '  -2.29%  '
$ws.Range('E35').Value = '  -2.27%  '
$ws.Range('D36').Value = '1.138.42'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('D37').Value = '''0.869'
$ws.Range('E37').Value = '  -4.13%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').Value = '''0.521'
$ws.Range('E39').Value = '  -3.69%  '
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('D41').Value = '''98.55'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('E44').Value = '  -4.84%  '
$ws.Range('D45').Value = '0.0₆0115'
$ws.Range('E45').Value = '  -1.95%  '
$ws.Range('D46').Value = '''55.16'
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').Value = '''0.415'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '''7.54'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('E51').Value = '  -0.02%  '
